# Minor text corrections on two slides.
$p = $ppt.ActivePresentation

# --- Slide 2: "Barista provides" bullet list -----------------------------
# Fix "client use the same SQL they used before" -> "clients use the same
# SQL they used before" (the sentence was originally split across two runs,
# "c" + "lient use the same SQL they used before"; merge/fix into one run).
$slide2 = $p.Slides.Item(2)
$tr2 = $slide2.Shapes.Item("Content Placeholder 2").TextFrame.TextRange

$para2 = $tr2.Paragraphs(8)
$text2 = $para2.Text.TrimEnd("`r")
if ($text2 -ne "client use the same SQL they used before") {
    throw "Slide 2 paragraph 8 text did not match expected content: [$text2]"
}

$start2 = $para2.Start
# Update the trailing run ("lient use the same SQL they used before") first
# so the merged run inherits its character formatting (dirty/smtClean), then
# delete the leading single-character run ("c") that is now redundant.
$tr2.Characters($start2 + 1, $text2.Length - 1).Text = "clients use the same SQL they used before"
$tr2.Characters($start2, 1).Text = ""

# --- Slide 5: "Design Choice: Enforcing Ordering" bullet list ------------
# Fix "... the transactions were submitted" -> "... the transactions are
# submitted", splitting the final run into "are " and "submitted".
$slide5 = $p.Slides.Item(5)
$tr5 = $slide5.Shapes.Item("Content Placeholder 2").TextFrame.TextRange

$para5 = $tr5.Paragraphs(4)
$text5 = $para5.Text.TrimEnd("`r")
$expected5 = "the commit order can be different from the order in which the transactions were submitted"
if ($text5 -ne $expected5) {
    throw "Slide 5 paragraph 4 text did not match expected content: [$text5]"
}

$start5 = $para5.Start
$offset5 = $text5.IndexOf("were submitted")
$tr5.Characters($start5 + $offset5, 5).Text = "are "
